$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1-3 header label corrections ---
$ws.Range("A1").Value = "TSR_RMIE"
$ws.Range("B1").Value = "TSR_RMIE"
$ws.Range("L1").Value = "ACBR_iBNO"
$ws.Range("A2").Value = "BTS"
$ws.Range("F2").Value = "1;w1,w1"
$ws.Range("G2").Value = "1;w1,w1"
$ws.Range("H2").Value = "1;w1,w1"
$ws.Range("I2").Value = "1;w1,w1"
$ws.Range("J2").Value = "BTS_TWW"
$ws.Range("A3").Value = "1;w1,w1"
$ws.Range("B3").Value = "TSR_RM"
$ws.Range("C3").Value = "TSR_RM"
$ws.Range("D3").Value = "1;w1,w1"
$ws.Range("F3").Value = "ACBR_BNO"
$ws.Range("G3").Value = "BTS"
$ws.Range("H3").Value = "ACBR_BNO"
$ws.Range("J3").Value = "TSR_TWW"

# --- Clear now-unused header cells (K2:M2, K3:M3) ---
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()

# --- Updated simulation result data (columns I and J) ---
$ws.Range("I6").Value = 5445
$ws.Range("J6").Value = 5133.75
$ws.Range("I7").Value = 4530
$ws.Range("J7").Value = 5133.75
$ws.Range("I8").Value = 3403.75
$ws.Range("J8").Value = 3233.75
$ws.Range("I9").Value = 5955
$ws.Range("J9").Value = 6080
$ws.Range("I10").Value = 2237.5
$ws.Range("J10").Value = 2401.25
$ws.Range("I11").Value = 4128.75
$ws.Range("J11").Value = 4227.5
$ws.Range("I12").Value = 4317.5
$ws.Range("J12").Value = 4315
$ws.Range("I13").Value = 7110
$ws.Range("J13").Value = 7157.5
$ws.Range("I14").Value = 3181.25
$ws.Range("J14").Value = 3263.75
$ws.Range("I15").Value = 4231.25
$ws.Range("J15").Value = 4325
$ws.Range("I16").Value = 1908.75
$ws.Range("J16").Value = 1897.5
$ws.Range("I17").Value = 3161.25
$ws.Range("J17").Value = 3161.25
$ws.Range("I18").Value = 2051.25
$ws.Range("J18").Value = 2032.5
$ws.Range("I19").Value = 2137.5
$ws.Range("J19").Value = 2076.25
$ws.Range("I21").Value = 5692.5
$ws.Range("J21").Value = 5695
$ws.Range("I22").Value = 4838.75
$ws.Range("J22").Value = 5377.5
$ws.Range("I23").Value = 3307.5
$ws.Range("J23").Value = 3311.25
$ws.Range("I24").Value = 6273.75
$ws.Range("J24").Value = 6663.75
$ws.Range("I25").Value = 2673.75
$ws.Range("J25").Value = 2978.75
$ws.Range("I26").Value = 4611.25
$ws.Range("J26").Value = 5075
$ws.Range("I27").Value = 4277.5
$ws.Range("J27").Value = 4090
$ws.Range("I28").Value = 6921.25
$ws.Range("J28").Value = 6966.25
$ws.Range("I29").Value = 3515
$ws.Range("J29").Value = 3505
$ws.Range("I30").Value = 3812.5
$ws.Range("J30").Value = 3885
$ws.Range("I31").Value = 1873.75
$ws.Range("J31").Value = 2006.25
$ws.Range("I32").Value = 2726.25
$ws.Range("J32").Value = 2817.5
$ws.Range("I33").Value = 2192.5
$ws.Range("J33").Value = 2237.5
$ws.Range("I34").Value = 2263.75
$ws.Range("J34").Value = 2245
